$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.258.75"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.329.55"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.58"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "2.324.35"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.62"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.67"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "2.744.33"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "60.246.53"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "2.327.45"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.10"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.90"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.17"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.93"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.38"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +10.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.36"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.06"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.383"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.04"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "323.52"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.10"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.95"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0949"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.41"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.566"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0498"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "0.0₆0219"
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.93"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
